# Update to test case excel file
#
# - "fundamentals" sheet: tester-name header cell and six new rows of
#   test-case data (word / letter-contained / expected result).
# - Selection moves to the "fundamentals" sheet (cell F9) and the
#   "manipulating" sheet's remembered selection moves to F6.

$wb = $excel.ActiveWorkbook

$fundamentals = $wb.Worksheets.Item("fundamentals")
$manipulating = $wb.Worksheets.Item("manipulating")

# Header / tester name cell.
$fundamentals.Range("A1").Value = "Manav"

# New test-case rows (2 columns of text contain column E duplicating the
# letter, column E = PASS, column F = explanatory comment).
$fundamentals.Cells.Item(3, 1).Value = 1
$fundamentals.Cells.Item(3, 2).Value = "Manav"
$fundamentals.Cells.Item(3, 3).Value = "n"
$fundamentals.Cells.Item(3, 4).Value = "n"
$fundamentals.Cells.Item(3, 5).Value = "PASS"
$fundamentals.Cells.Item(3, 6).Value = "The function is working as expected."

$fundamentals.Cells.Item(4, 1).Value = 2
$fundamentals.Cells.Item(4, 2).Value = "Seneca"
$fundamentals.Cells.Item(4, 3).Value = "c"
$fundamentals.Cells.Item(4, 4).Value = "c"
$fundamentals.Cells.Item(4, 5).Value = "PASS"
$fundamentals.Cells.Item(4, 6).Value = "The function is working as expected."

$fundamentals.Cells.Item(5, 1).Value = 3
$fundamentals.Cells.Item(5, 2).Value = "Toronto"
$fundamentals.Cells.Item(5, 3).Value = "o"
$fundamentals.Cells.Item(5, 4).Value = "o"
$fundamentals.Cells.Item(5, 5).Value = "PASS"
$fundamentals.Cells.Item(5, 6).Value = "The function is working as expected."

$fundamentals.Cells.Item(6, 1).Value = 4
$fundamentals.Cells.Item(6, 2).Value = "Markham"
$fundamentals.Cells.Item(6, 3).Value = "a"
$fundamentals.Cells.Item(6, 4).Value = "a"
$fundamentals.Cells.Item(6, 5).Value = "PASS"
$fundamentals.Cells.Item(6, 6).Value = "The function is working as expected."

$fundamentals.Cells.Item(7, 1).Value = 5
$fundamentals.Cells.Item(7, 2).Value = "Canada"
$fundamentals.Cells.Item(7, 3).Value = "a"
$fundamentals.Cells.Item(7, 4).Value = "a"
$fundamentals.Cells.Item(7, 5).Value = "PASS"
$fundamentals.Cells.Item(7, 6).Value = "The function is working as expected."

$fundamentals.Cells.Item(8, 1).Value = 6
$fundamentals.Cells.Item(8, 2).Value = "Ontario"
$fundamentals.Cells.Item(8, 3).Value = "O"
$fundamentals.Cells.Item(8, 4).Value = "O"
$fundamentals.Cells.Item(8, 5).Value = "PASS"
$fundamentals.Cells.Item(8, 6).Value = "The function is working as expected."

# Leave "manipulating" selected at F6 first, then select "fundamentals" at
# F9 last so it ends up as the active / tabSelected sheet.
[void]$manipulating.Range("F6").Select()
[void]$fundamentals.Range("F9").Select()
